$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainText($rangeAddr, $val) {
    $ws.Range($rangeAddr).Value = $val
}

function Set-NumericLookingText($rangeAddr, $val) {
    $cell = $ws.Range($rangeAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-PlainText "D2" '38.831.79'
Set-PlainText "D3" '2.091.23'
Set-PlainText "E3" '  +2.25%  '
Set-PlainText "E4" '  +0.02%  '
Set-NumericLookingText "D5" '228.27'
Set-PlainText "E6" '  +0.63%  '
Set-NumericLookingText "D7" '60.52'
Set-PlainText "E7" '  +0.47%  '
Set-PlainText "E8" '  +0.00%  '
Set-PlainText "E9" '  +2.21%  '
Set-NumericLookingText "D10" '0.0837'
Set-PlainText "E10" '  +0.38%  '
Set-NumericLookingText "D11" '0.103'
Set-PlainText "E11" '  -0.65%  '
Set-PlainText "D12" '2.402.29'
Set-PlainText "E12" '  +2.31%  '
Set-NumericLookingText "D13" '14.96'
Set-PlainText "E13" '  +3.90%  '
Set-NumericLookingText "D14" '21.96'
Set-PlainText "E14" '  +2.35%  '
Set-PlainText "E15" '  +4.08%  '
Set-PlainText "E16" '  -0.80%  '
Set-PlainText "D17" '2.092.34'
Set-PlainText "E17" '  +2.27%  '
Set-PlainText "D18" '38.763.49'
Set-PlainText "E18" '  +2.70%  '
Set-PlainText "E19" '  +3.25%  '
Set-PlainText "E20" '  +2.10%  '
Set-PlainText "D21" '0.0₃0836'
Set-PlainText "E21" '  +1.06%  '
Set-NumericLookingText "D22" '227.26'
Set-PlainText "E22" '  +2.12%  '
Set-PlainText "E23" '  -0.43%  '
Set-NumericLookingText "D24" '2.38'
Set-PlainText "E24" '  -0.41%  '
Set-PlainText "E25" '  +3.06%  '
Set-NumericLookingText "D26" '170.88'
Set-PlainText "E26" '  +1.08%  '
Set-PlainText "E27" '  +1.94%  '
Set-NumericLookingText "D28" '0.141'
Set-PlainText "E28" '  +9.61%  '
Set-NumericLookingText "D29" '1.48'
Set-PlainText "E29" '  +14.82%  '
Set-NumericLookingText "D30" '19.18'
Set-PlainText "E30" '  +2.21%  '
Set-PlainText "E31" '  +0.82%  '
Set-NumericLookingText "D32" '2.37'
Set-PlainText "E32" '  +5.39%  '
Set-PlainText "E33" '  +2.77%  '
Set-PlainText "E34" '  +4.30%  '
Set-PlainText "E35" '  +1.42%  '
Set-NumericLookingText "D36" '6.46'
Set-PlainText "E36" '  -0.80%  '
Set-PlainText "E37" '  +1.65%  '
Set-NumericLookingText "D38" '3.59'
Set-PlainText "E38" '  +3.64%  '
Set-PlainText "E39" '  -0.01%  '
Set-PlainText "E40" '  -0.19%  '
Set-PlainText "D41" '1.541.27'
Set-PlainText "E41" '  +0.59%  '
Set-PlainText "B42" 'VeChain'
Set-PlainText "C42" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-NumericLookingText "D42" '0.0225'
Set-PlainText "E42" '  +4.28%  '
Set-PlainText "B43" 'Aave'
Set-PlainText "C43" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-NumericLookingText "D43" '100.90'
Set-PlainText "E43" '  +3.18%  '
Set-PlainText "E44" '  -0.79%  '
Set-NumericLookingText "D45" '0.0921'
Set-PlainText "E45" '  +3.44%  '
Set-PlainText "E46" '  +9.00%  '
Set-NumericLookingText "D47" '1.12'
Set-PlainText "E47" '  +1.35%  '
Set-NumericLookingText "D48" '4.13'
Set-PlainText "E48" '  -0.38%  '
Set-PlainText "E49" '  +2.91%  '
Set-PlainText "E50" '  +1.11%  '
Set-PlainText "D51" '2.289.20'
Set-PlainText "E51" '  +2.42%  '
